# ---------------------------------------------------------------------------
# "inserido o crud do cadastro de lancamento"
#
# Inserts a new "lancamento" worksheet (CRUD-doc sheet, same layout as the
# existing "servico" sheet) right before "Planilha2", makes it the active
# sheet/tab, and tweaks the view-state of a couple of other sheets
# (tabSelected / scroll-selection / a row height) the way the original
# author's Excel session left them.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the "lancamento" sheet immediately before "Planilha2"
# ---------------------------------------------------------------------------
$planilha2 = $wb.Worksheets.Item("Planilha2")
$ws = $wb.Worksheets.Add($planilha2)
$ws.Name = "lancamento"

# ---------------------------------------------------------------------------
# 2) Header row (1) -- field names
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "descricao"
$ws.Range("C1").Value = "tipo"
$ws.Range("D1").Value = "valor"
$ws.Range("E1").Value = "dataCadastro"
$ws.Range("F1").Value = "dataModificacao"
$ws.Range("G1").Value = "dataVencimento"
$ws.Range("H1").Value = "pago"

# ---------------------------------------------------------------------------
# 3) Sample data row (2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Financeiro"
$ws.Range("C2").Value = "Entrada"

# D2 = text "150.98" (not a number) -- force text without leaving a
# quote-prefix style behind: stage it on a scratch cell, copy the *value*
# across, then wipe the scratch cell (content + format) completely.
$ws.Range("ZZ1").Value = "'150.98"
$ws.Range("ZZ1").Copy()
$ws.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

# E2/F2/G2 share the same datetime serial + the existing "m/d/yyyy h:mm"
# style already used elsewhere in the workbook (copy the format so it
# reuses the same style record instead of minting a new numFmt).
$servico = $wb.Worksheets.Item("servico")
$servico.Range("F2").Copy()
$ws.Range("E2:G2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("E2").Value = 43986.545439814814
$ws.Range("F2").Value = 43986.545439814814
$ws.Range("G2").Value = 43986.545439814814
$ws.Range("H2").Value = 1

# ---------------------------------------------------------------------------
# 4) Second header row (4) -- bold "Campo / Valor / Atributos PHP / ..."
# ---------------------------------------------------------------------------
$ws.Range("A4:H4").Font.Bold = $true
$ws.Range("A4").Value = "Campo"
$ws.Range("B4").Value = "Valor"
$ws.Range("C4").Value = "Atributos PHP"
$ws.Range("D4").Value = "Get PHP"
$ws.Range("E4").Value = "Set PHP"
$ws.Range("F4").Value = "Name Form Html"
$ws.Range("G4").Value = "Observação html"
$ws.Range("H4").Value = "Set Controle"

# ---------------------------------------------------------------------------
# 5) CRUD snippet rows (5-12), one per field, mirroring the "servico" sheet
# ---------------------------------------------------------------------------
$ws.Range("D5:E12").WrapText = $true

$fields = @(
    @{ Row = 5;  Name = "id";              BVal = 2;                      BIsNum = $true;  Height = 86.4  },
    @{ Row = 6;  Name = "descricao";       BVal = "Financeiro";           BIsNum = $false; Height = 100.8 },
    @{ Row = 7;  Name = "tipo";            BVal = "Entrada";              BIsNum = $false; Height = 86.4  },
    @{ Row = 8;  Name = "valor";           BVal = "150.98";               BIsNum = $false; Height = 115.2 },
    @{ Row = 9;  Name = "dataCadastro";    BVal = 43986.545439814814;     BIsNum = $true;  Height = 144   },
    @{ Row = 10; Name = "dataModificacao"; BVal = 43986.545439814814;     BIsNum = $true;  Height = 129.6 },
    @{ Row = 11; Name = "dataVencimento";  BVal = 43986.545439814814;     BIsNum = $true;  Height = 129.6 },
    @{ Row = 12; Name = "pago";            BVal = 1;                      BIsNum = $true;  Height = 115.2 }
)

foreach ($fld in $fields) {
    $row = $fld.Row

    $ws.Rows.Item($row).RowHeight = $fld.Height

    $ws.Range("A$row").Value = $fld.Name
    if ($fld.BIsNum) {
        $ws.Range("B$row").Value = $fld.BVal
    } else {
        $ws.Range("B$row").Value = $fld.BVal
    }

    $ws.Range("C$row").Formula = "=""$""&TRIM(A$row)&"","""

    $getF = "=""public function get""&TRIM(`$A$row)&""(){" + "`n`t`treturn `$this->""&TRIM(`$A$row)&"";" + "`n`t}"""
    $ws.Range("D$row").Formula = $getF

    $setF = "=""public function set""&TRIM(`$A$row)&""(`$""&TRIM(A$row)&""){" + "`n`t`t`$this->""&TRIM(`$A$row)&"" = `$""&TRIM(`$A$row)&"";" + "`n`t}"""
    $ws.Range("E$row").Formula = $setF

    $ws.Range("F$row").Formula = "=""name=""""""&A$row&""""""""""

    $ws.Range("H$row").Formula = "=""`$lancamento->set""&`$A$row&""(`$_POST['""&`$A$row&""']);"""
}

# B8 ("valor" sample) must stay plain text "150.98", same trick as D2 above.
$ws.Range("ZZ1").Value = "'150.98"
$ws.Range("ZZ1").Copy()
$ws.Range("B8").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 6) Column widths (best-fit-ish, closest the host lets us set)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 14.1666666
$ws.Columns.Item(2).ColumnWidth = 11.1666666
$ws.Columns.Item(3).ColumnWidth = 15.6
$ws.Columns.Item(4).ColumnWidth = 11.5
$ws.Columns.Item(5).ColumnWidth = 16.9333333
$ws.Columns.Item(6).ColumnWidth = 21.5
$ws.Columns.Item(7).ColumnWidth = 14.8333333
$ws.Columns.Item(8).ColumnWidth = 50.0533333

# ---------------------------------------------------------------------------
# 7) Sheet view / selection + page margins matching the rest of the workbook
# ---------------------------------------------------------------------------
$ws.PageSetup.LeftMargin = 36.850393728
$ws.PageSetup.RightMargin = 36.850393728
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995
$ws.PageSetup.HeaderMargin = 22.67716464
$ws.PageSetup.FooterMargin = 22.67716464

$ws.Activate()
$ws.Range("F10").Select()

# ---------------------------------------------------------------------------
# 8) "usuario" sheet -- it's no longer the active tab, and row 6 shrank
# ---------------------------------------------------------------------------
$usuario = $wb.Worksheets.Item("usuario")
$usuario.Activate()
$usuario.Range("E5:E14").Select()
$usuario.Rows.Item(6).RowHeight = 72

# ---------------------------------------------------------------------------
# 9) "servico" sheet -- scrolled/selected differently now
# ---------------------------------------------------------------------------
$servico.Activate()
$servico.Range("A4:H12").Select()

# ---------------------------------------------------------------------------
# 10) Leave "lancamento" as the active sheet/tab, like the saved workbook
# ---------------------------------------------------------------------------
$ws.Activate()
